# "Generate Report for Handoff"
#
# The 6afa548a-*.md file was already handed back and is in sync with en-US.
# The b30a2485-*.md file has now gone through a new handoff cycle (new
# handoff timestamps), so its row moves to the top of each report and the
# other file's status flips to "Ready for handoff". This script rewrites
# the three report sheets (Overview, zh-cn, de-de) accordingly, swapping the
# two data rows and updating the handoff timestamps / status text, while
# fixing up each hyperlink's visible display text to match its new cell.

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay($ws, [string]$addr, [string]$text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

$docMd = "6afa548a-6a5f-4775-b2ab-9ec8106b51ab.md"
$depMd = "b30a2485-648f-495c-928b-a6f7cff6434b.md"

$statusInSync = "Handed back: in sync with en-US"
$statusReady = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("A2").Value = $depMd
$ovw.Range("B2").Value = $statusInSync
$ovw.Range("C2").Value = $statusInSync

$ovw.Range("A3").Value = $docMd
$ovw.Range("B3").Value = $statusReady
$ovw.Range("C3").Value = $statusReady

Set-HyperlinkDisplay $ovw '$A$2' $depMd
Set-HyperlinkDisplay $ovw '$A$3' $docMd

# ---------------------------------------------------------------------
# Per-locale sheets (zh-cn / de-de)
# ---------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; XlfDoc = "6afa548a-6a5f-4775-b2ab-9ec8106b51ab.ea15ff70dcf6c2548779b068e0543e78b4e61bfc.zh-cn.xlf"; XlfDep = "b30a2485-648f-495c-928b-a6f7cff6434b.780fa0f014ca4f8a334bf35157504a83a44af9b9.zh-cn.xlf"; HandoffTime = "2016-02-23 08:17:24"; HandbackTime = "2016-02-23 08:16:00" },
    @{ Sheet = "de-de"; XlfDoc = "6afa548a-6a5f-4775-b2ab-9ec8106b51ab.ea15ff70dcf6c2548779b068e0543e78b4e61bfc.de-de.xlf"; XlfDep = "b30a2485-648f-495c-928b-a6f7cff6434b.780fa0f014ca4f8a334bf35157504a83a44af9b9.de-de.xlf"; HandoffTime = "2016-02-23 08:17:40"; HandbackTime = "2016-02-23 08:16:27" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Row 2 -> dependency file (b30a2485), freshly handed off
    $ws.Range("A2").Value = $depMd
    $ws.Range("B2").Value = $statusInSync
    $ws.Range("C2").Value = $loc.XlfDep
    $ws.Range("D2").Value = $loc.HandoffTime
    $ws.Range("E2").Value = $depMd
    $ws.Range("F2").Value = $loc.XlfDep
    $ws.Range("G2").Value = $loc.HandbackTime
    $ws.Range("H2").Value = "Include"

    # Row 3 -> main doc file (6afa548a), now ready for handoff again
    $ws.Range("A3").Value = $docMd
    $ws.Range("B3").Value = $statusReady
    $ws.Range("C3").Value = $loc.XlfDoc
    $ws.Range("D3").Value = $loc.HandoffTime
    $ws.Range("E3").Value = $docMd
    $ws.Range("F3").Value = $loc.XlfDoc
    $ws.Range("G3").Value = $loc.HandbackTime
    $ws.Range("H3").Value = "Include"

    Set-HyperlinkDisplay $ws '$A$2' $depMd
    Set-HyperlinkDisplay $ws '$C$2' $loc.XlfDep
    Set-HyperlinkDisplay $ws '$E$2' $depMd
    Set-HyperlinkDisplay $ws '$F$2' $loc.XlfDep

    Set-HyperlinkDisplay $ws '$A$3' $docMd
    Set-HyperlinkDisplay $ws '$C$3' $loc.XlfDoc
    Set-HyperlinkDisplay $ws '$E$3' $docMd
    Set-HyperlinkDisplay $ws '$F$3' $loc.XlfDoc
}
